$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 384.3
$ws.Range("I12").Value = 204.77777
$ws.Range("K12").Value = 204.77777
$ws.Range("M12").Value = -34.77777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 53942
$ws.Range("J93").Value = 53942
$ws.Range("L93").Value = 53942
$ws.Range("N93").Value = -58934

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 125000690
$ws.Range("I125").Value = 166666800
$ws.Range("K125").Value = 1500001200
$ws.Range("M125").Value = -1499998740

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 924
$ws.Range("I132").Value = 924
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2772
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -242
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2004487.9
$ws.Range("I138").Value = 2212.2856
$ws.Range("J138").Value = 3454411.5
$ws.Range("K138").Value = 6636.8568
$ws.Range("L138").Value = 10363234.5
$ws.Range("M138").Value = -1496.8568
$ws.Range("N138").Value = -10373514.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4173466.5
$ws.Range("I32").Value = 4352965
$ws.Range("K32").Value = 4352965
$ws.Range("M32").Value = -4352678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2038.8
$ws.Range("I35").Value = 1673.625
$ws.Range("J35").Value = 3499.5
$ws.Range("K35").Value = 1673.625
$ws.Range("L35").Value = 3499.5
$ws.Range("M35").Value = -1267.625
$ws.Range("N35").Value = -4311.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3622.5757
$ws.Range("I61").Value = 2745.4211
$ws.Range("J61").Value = 4813
$ws.Range("K61").Value = 2745.4211
$ws.Range("L61").Value = 4813
$ws.Range("M61").Value = -2533.4211
$ws.Range("N61").Value = -5237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 676.1111
$ws.Range("I88").Value = 635.2
$ws.Range("J88").Value = 727.25
$ws.Range("K88").Value = 635.2
$ws.Range("L88").Value = 727.25
$ws.Range("M88").Value = -229.2
$ws.Range("N88").Value = -1539.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 676.1111
$ws.Range("I91").Value = 635.2
$ws.Range("J91").Value = 727.25
$ws.Range("K91").Value = 635.2
$ws.Range("L91").Value = 727.25
$ws.Range("M91").Value = 768.8
$ws.Range("N91").Value = -3535.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 7409980
$ws.Range("I102").Value = 14288180
$ws.Range("K102").Value = 14288180
$ws.Range("M102").Value = -14286558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3610.3428
$ws.Range("I122").Value = 2042.091
$ws.Range("K122").Value = 6126.272999999999
$ws.Range("M122").Value = -3676.272999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3371.3572
$ws.Range("I132").Value = 1395.1724
$ws.Range("K132").Value = 4185.5172
$ws.Range("M132").Value = -1655.5172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3622.5757
$ws.Range("I136").Value = 2745.4211
$ws.Range("J136").Value = 4813
$ws.Range("K136").Value = 8236.263300000001
$ws.Range("L136").Value = 14439
$ws.Range("M136").Value = -5686.263300000001
$ws.Range("N136").Value = -19539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 51476.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 51476.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 51476.668
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -61756.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 41500
$ws.Range("J69").Value = 41500
$ws.Range("L69").Value = 41500
$ws.Range("N69").Value = -43122

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 41500
$ws.Range("J72").Value = 41500
$ws.Range("L72").Value = 124500
$ws.Range("N72").Value = -132612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4796.9395
$ws.Range("I31").Value = 2374.5293
$ws.Range("K31").Value = 2374.5293
$ws.Range("M31").Value = -2079.5293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4796.9395
$ws.Range("I34").Value = 2374.5293
$ws.Range("K34").Value = 2374.5293
$ws.Range("M34").Value = -2172.5293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4300.483
$ws.Range("I62").Value = 4025.4736
$ws.Range("J62").Value = 4823
$ws.Range("K62").Value = 4025.4736
$ws.Range("L62").Value = 4823
$ws.Range("M62").Value = -3401.4736
$ws.Range("N62").Value = -6071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4300.483
$ws.Range("I65").Value = 4025.4736
$ws.Range("J65").Value = 4823
$ws.Range("K65").Value = 20127.368
$ws.Range("L65").Value = 24115
$ws.Range("M65").Value = -17007.368
$ws.Range("N65").Value = -30355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7937694.5
$ws.Range("I105").Value = 8929282
$ws.Range("J105").Value = 4995
$ws.Range("K105").Value = 8929282
$ws.Range("L105").Value = 4995
$ws.Range("M105").Value = -8927535
$ws.Range("N105").Value = -8489

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2126.8948
$ws.Range("I107").Value = 716.25
$ws.Range("J107").Value = 3152.818
$ws.Range("K107").Value = 716.25
$ws.Range("L107").Value = 3152.818
$ws.Range("M107").Value = 1203.75
$ws.Range("N107").Value = -6992.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3552.2
$ws.Range("I134").Value = 1850.4348
$ws.Range("K134").Value = 5551.3044
$ws.Range("M134").Value = -3016.3044

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 658.619
$ws.Range("I5").Value = 591.55
$ws.Range("K5").Value = 1774.65
$ws.Range("M5").Value = -1662.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25004370
$ws.Range("I68").Value = 66667548
$ws.Range("J68").Value = 6464.4
$ws.Range("K68").Value = 200002644
$ws.Range("L68").Value = 19393.2
$ws.Range("M68").Value = -200001833
$ws.Range("N68").Value = -21015.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 25004370
$ws.Range("I71").Value = 66667548
$ws.Range("J71").Value = 6464.4
$ws.Range("K71").Value = 600007932
$ws.Range("L71").Value = 58179.6
$ws.Range("M71").Value = -600003876
$ws.Range("N71").Value = -66291.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2328.3137
$ws.Range("I131").Value = 1775.4445
$ws.Range("J131").Value = 2446.7856
$ws.Range("K131").Value = 5326.333500000001
$ws.Range("L131").Value = 7340.3568
$ws.Range("M131").Value = -286.3335000000006
$ws.Range("N131").Value = -17420.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 7418.6523
$ws.Range("I132").Value = 2629
$ws.Range("J132").Value = 12643.728
$ws.Range("K132").Value = 23661
$ws.Range("L132").Value = 113793.552
$ws.Range("M132").Value = -21131
$ws.Range("N132").Value = -118853.552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 658.619
$ws.Range("I135").Value = 591.55
$ws.Range("K135").Value = 5323.95
$ws.Range("M135").Value = -2788.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 97940.14
$ws.Range("I137").Value = 101554.3
$ws.Range("J137").Value = 94654.55
$ws.Range("K137").Value = 304662.9
$ws.Range("L137").Value = 283963.65
$ws.Range("M137").Value = -299562.9
$ws.Range("N137").Value = -294163.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7937.5884
$ws.Range("I70").Value = 5332.6665
$ws.Range("J70").Value = 8495.786
$ws.Range("K70").Value = 5332.6665
$ws.Range("L70").Value = 8495.786
$ws.Range("M70").Value = -5062.6665
$ws.Range("N70").Value = -9035.786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7937.5884
$ws.Range("I73").Value = 5332.6665
$ws.Range("J73").Value = 8495.786
$ws.Range("K73").Value = 5332.6665
$ws.Range("L73").Value = 8495.786
$ws.Range("M73").Value = -4396.6665
$ws.Range("N73").Value = -10367.786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5240
$ws.Range("I80").Value = 1246.5
$ws.Range("J80").Value = 6381
$ws.Range("K80").Value = 1246.5
$ws.Range("L80").Value = 6381
$ws.Range("M80").Value = -248.5
$ws.Range("N80").Value = -8377

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5240
$ws.Range("I83").Value = 1246.5
$ws.Range("J83").Value = 6381
$ws.Range("K83").Value = 6232.5
$ws.Range("L83").Value = 31905
$ws.Range("M83").Value = -1240.5
$ws.Range("N83").Value = -41889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5589.6216
$ws.Range("I113").Value = 2161.8462
$ws.Range("J113").Value = 7446.3335
$ws.Range("K113").Value = 2161.8462
$ws.Range("L113").Value = 7446.3335
$ws.Range("M113").Value = 8.153800000000047
$ws.Range("N113").Value = -11786.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1940.4706
$ws.Range("I132").Value = 1622.683
$ws.Range("K132").Value = 4868.049
$ws.Range("M132").Value = -2338.049

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4632419
$ws.Range("I46").Value = 2321.3572
$ws.Range("J46").Value = 11114555
$ws.Range("K46").Value = 2321.3572
$ws.Range("L46").Value = 11114555
$ws.Range("M46").Value = -2133.3572
$ws.Range("N46").Value = -11114931

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4938.3105
$ws.Range("J61").Value = 7440.1665
$ws.Range("L61").Value = 7440.1665
$ws.Range("N61").Value = -7844.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 835832
$ws.Range("I82").Value = 2000720.4
$ws.Range("J82").Value = 3768.8572
$ws.Range("K82").Value = 2000720.4
$ws.Range("L82").Value = 3768.8572
$ws.Range("M82").Value = -2000359.4
$ws.Range("N82").Value = -4490.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 835832
$ws.Range("I85").Value = 2000720.4
$ws.Range("J85").Value = 3768.8572
$ws.Range("K85").Value = 2000720.4
$ws.Range("L85").Value = 3768.8572
$ws.Range("M85").Value = -1999472.4
$ws.Range("N85").Value = -6264.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4938.3105
$ws.Range("J113").Value = 7440.1665
$ws.Range("L113").Value = 7440.1665
$ws.Range("N113").Value = -11780.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 30000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -31872

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 90000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -99360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10288625
$ws.Range("I122").Value = 13264949
$ws.Range("J122").Value = 6781
$ws.Range("K122").Value = 39794847
$ws.Range("L122").Value = 20343
$ws.Range("M122").Value = -39792397
$ws.Range("N122").Value = -25243

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4836.161
$ws.Range("I132").Value = 4612.68
$ws.Range("K132").Value = 13838.04
$ws.Range("M132").Value = -11308.04

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19059620
$ws.Range("I136").Value = 33334326
$ws.Range("J136").Value = 440439.75
$ws.Range("K136").Value = 100002978
$ws.Range("L136").Value = 1321319.25
$ws.Range("M136").Value = -100000428
$ws.Range("N136").Value = -1326419.25
